# edit.ps1 -- reproduces the OOXML diff against HarishKumarBR_5Years_DevOps.docx
#
# Strategy: every change below is expressed as a *whole paragraph* replacement.
# We locate a paragraph via unique anchor text (Find on rendered text, which is
# robust to runs being split), expand the hit to the full paragraph (including
# its end-of-paragraph mark) with Range.Expand(wdParagraph), and then push the
# exact target OOXML for that paragraph back in with Range.InsertXML wrapped in
# a minimal single-part WordProcessingML package. This preserves pPr/numbering
# that a naive Range.Text edit would otherwise clobber, and lets us add/drop
# elements (like <w:lastRenderedPageBreak/>) that have no dedicated object
# model property.

$wdParagraph = 4

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the paragraph containing $anchorText (searching the whole document,
# or from $searchFrom onward when given) and overwrite that whole paragraph
# with $newParaXml (a single <w:p>...</w:p> fragment). Returns the (reseated)
# range that was written, positioned at the start of the replaced paragraph,
# or $null if the anchor wasn't found.
function Set-ParagraphXmlByAnchor([string]$anchorText, [string]$newParaXml, $searchFrom) {
    if ($searchFrom -eq $null) {
        $searchFrom = 0
    }
    $hit = $d.Range($searchFrom, $d.Content.End)
    $found = $hit.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "ANCHOR NOT FOUND:" $anchorText
        return $null
    }
    $paraRng = $d.Range($hit.Start, $hit.Start)
    $paraRng.Expand($wdParagraph)
    $paraRng.InsertXML((New-PkgXml $newParaXml))
    return $paraRng
}

# ---------------------------------------------------------------------------
# 1) Duplicate the empty "section separator" paragraph (bold/underlined,
#    tab stop at 1680, 360-line spacing) that sits right before the first
#    "Roles and Responsibilities:" heading (the SQE / ivpn2,MX,connect
#    project block), inserting a second copy of it just above that heading.
# ---------------------------------------------------------------------------

$sepParaXml = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1680"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>'

# Anchor on the unique preceding sentence, then find the (unique, first)
# "Roles and Responsibilities:" after it.
$sqeHit = $d.Range(0, $d.Content.End)
$sqeFound = $sqeHit.Find.Execute("ivpn2,MX,connect", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $sqeFound) {
    Write-Host "ANCHOR NOT FOUND: ivpn2,MX,connect"
} else {
    $rolesHit = $d.Range($sqeHit.End, $d.Content.End)
    $rolesFound = $rolesHit.Find.Execute("Roles and Responsibilities:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rolesFound) {
        Write-Host "ANCHOR NOT FOUND: Roles and Responsibilities: (after SQE paragraph)"
    } else {
        # The existing empty separator paragraph ends exactly one character
        # (its paragraph mark) before the "Roles and Responsibilities:" text.
        $existingSep = $d.Range($rolesHit.Start - 1, $rolesHit.Start - 1)
        $existingSep.Expand($wdParagraph)
        # Replace that single separator paragraph with two copies of itself.
        $existingSep.InsertXML((New-PkgXml ($sepParaXml + $sepParaXml)))
    }
}

# ---------------------------------------------------------------------------
# 2) Split the "Strong experience in Kubernetes..." run in two, inserting a
#    <w:lastRenderedPageBreak/> right before the final word "cluster.".
# ---------------------------------------------------------------------------

$strongXml = '<w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:suppressAutoHyphens/><w:spacing w:line="360" w:lineRule="auto"/><w:contextualSpacing/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Strong experience in Kubernetes cluster setup and deploying applications in Kubernetes </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>cluster.</w:t></w:r></w:p>'

Set-ParagraphXmlByAnchor "Strong experience in Kubernetes cluster setup and deploying applications in Kubernetes cluster." $strongXml $null | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the <w:lastRenderedPageBreak/> that used to sit on the "Implemented
#    CI/CD using " run (the page break moved earlier, into the Kubernetes
#    paragraph above).
# ---------------------------------------------------------------------------

$cicdXml = '<w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:suppressAutoHyphens/><w:spacing w:line="360" w:lineRule="auto"/><w:contextualSpacing/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Implemented CI/CD using </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Jenkins,</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> D</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ocker &amp; Kubernetes.</w:t></w:r></w:p>'

Set-ParagraphXmlByAnchor "Implemented CI/CD using" $cicdXml $null | Out-Null

# ---------------------------------------------------------------------------
# 4) Add <w:lastRenderedPageBreak/> onto "Planning and documentation of the
#    requirements." (the page break moved here)...
# ---------------------------------------------------------------------------

$planningXml = '<w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:suppressAutoHyphens/><w:spacing w:line="360" w:lineRule="auto"/><w:contextualSpacing/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>Planning and documentation of the requirements.</w:t></w:r></w:p>'

Set-ParagraphXmlByAnchor "Planning and documentation of the requirements." $planningXml $null | Out-Null

# ---------------------------------------------------------------------------
# 5) ...and remove it from "Deploying the Jar, War, Ear files in WebLogic."
# ---------------------------------------------------------------------------

$weblogicXml = '<w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:suppressAutoHyphens/><w:spacing w:line="360" w:lineRule="auto"/><w:contextualSpacing/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Deploying the Jar, War, Ear files in WebLogic.</w:t></w:r></w:p>'

Set-ParagraphXmlByAnchor "Deploying the Jar, War, Ear files in WebLogic." $weblogicXml $null | Out-Null

# ---------------------------------------------------------------------------
# 6) The blank "List Paragraph" styled spacer right before "Certifications
#    and Achievements:" loses its pStyle and eastAsia font override, and its
#    left indent changes from 450 to 90 twips.
# ---------------------------------------------------------------------------

$spacerXml = '<w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:line="360" w:lineRule="auto"/><w:ind w:left="90"/><w:jc w:val="both"/><w:outlineLvl w:val="3"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr></w:p>'

# Anchor on the following unique heading and step back one paragraph to the
# (empty, text-less) spacer paragraph.
$certHit = $d.Range(0, $d.Content.End)
$certFound = $certHit.Find.Execute("Certifications and Achievements:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $certFound) {
    Write-Host "ANCHOR NOT FOUND: Certifications and Achievements:"
} else {
    $spacerRng = $d.Range($certHit.Start - 1, $certHit.Start - 1)
    $spacerRng.Expand($wdParagraph)
    $spacerRng.InsertXML((New-PkgXml $spacerXml))
}

# ---------------------------------------------------------------------------
# 7) word/numbering.xml also renumbers two orphaned VML picture-bullet shape
#    ids (_x0000_i1058 -> _x0000_i1034, _x0000_i1059 -> _x0000_i1035). Those
#    ids aren't attached to any list level (no <w:lvlPicBulletId> references
#    them) and there is no Word object-model surface (ListFormat/ListTemplate/
#    ListLevel.PictureBullet etc.) that reaches raw numPicBullet VML shape
#    ids, so this purely-cosmetic bookkeeping id isn't reproducible from
#    COM-interop automation; intentionally left as-is.
# ---------------------------------------------------------------------------

Write-Host "edit.ps1 completed"
